$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: human-readable labels that drive the A (hex length) and C (formatted macro line) formulas.
$ws.Range("B1").Value = 'CURVE.FORWARD.RATE'
$ws.Range("B2").Value = 'curve identifier,time 1,time 2,convention,trigger'
$ws.Range("B4").Value = 'xlObjectTools - Yield Curve'
$ws.Range("B7").Value = 'This function computes the forward rate for the corresponding dates'
$ws.Range("B8").Value = 'Curve identifier'
$ws.Range("B9").Value = 'The forward date'
$ws.Range("B10").Value = 'The maturity date'
$ws.Range("B11").Value = 'The identifier of the convention'
$ws.Range("B12").Value = 'Trigger for recalculation  '

# Column H: literal xlObjectTools macro-registration lines (static text, independent of the formulas).
$ws.Range("H1").Value = '            TempStrNoSize("\x13""COURBE.TAUX.FORWARD"),'
$ws.Range("H2").Value = '            TempStrNoSize("\x31""Identifiant de la courbe,t1,t2,convention,trigger"),'
$ws.Range("H4").Value = '            TempStrNoSize("\x1C""xlObjectTools - Yield Curve"),'
$ws.Range("H7").Value = '            TempStrNoSize("\x3B""Cette fonction extrait le taux forward de la courbe pointée"),'
$ws.Range("H8").Value = '            TempStrNoSize("\x11""Object Identifer"),'
$ws.Range("H9").Value = '            TempStrNoSize("\x2A""La date de positionnement de l''observation"),'
$ws.Range("H10").Value = '            TempStrNoSize("\x26""La date de terminaison du taux forward"),'
$ws.Range("H11").Value = '            TempStrNoSize("\x29""L''indentifiant de la convention de calcul"),'
$ws.Range("H12").Value = '            TempStrNoSize("\x17""Déclenche le recalcul  "));'

